$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Value)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.407.00'
Set-TextValue 'E2' '  -2.12%  '
Set-TextValue 'D3' '1.845.17'
Set-TextValue 'E3' '  -1.92%  '
Set-TextValue 'E4' '  -0.05%  '
Set-TextValue 'D5' '260.73'
Set-TextValue 'E5' '  -7.77%  '
Set-TextValue 'E6' '  +0.02%  '
Set-TextValue 'D7' '0.5264'
Set-TextValue 'E7' '  -0.09%  '
Set-TextValue 'D8' '0.3234'
Set-TextValue 'E8' '  -8.79%  '
Set-TextValue 'E9' '  -4.66%  '
Set-TextValue 'D10' '18.93'
Set-TextValue 'E10' '  -7.18%  '
Set-TextValue 'D11' '0.7716'
Set-TextValue 'E11' '  -5.92%  '
Set-TextValue 'D12' '0.07715'
Set-TextValue 'E12' '  -1.47%  '
Set-TextValue 'D13' '1.884.14'
Set-TextValue 'E13' '  -0.25%  '
Set-TextValue 'D14' '89.10'
Set-TextValue 'E14' '  -1.79%  '
Set-TextValue 'D15' '5.031'
Set-TextValue 'E15' '  -3.96%  '
Set-TextValue 'D16' '1.001'
Set-TextValue 'E16' '  -0.04%  '
Set-TextValue 'E17' '  -3.05%  '
Set-TextValue 'E18' '  +0.07%  '
Set-TextValue 'D19' '0.000007903'
Set-TextValue 'E19' '  -3.36%  '
Set-TextValue 'D20' '26.421.39'
Set-TextValue 'E20' '  -2.21%  '
Set-TextValue 'D21' '2.067.00'
Set-TextValue 'E21' '  -3.22%  '
Set-TextValue 'D22' '4.539'
Set-TextValue 'E22' '  -5.43%  '
Set-TextValue 'D23' '9.523'
Set-TextValue 'E23' '  -6.83%  '
Set-TextValue 'D24' '5.928'
Set-TextValue 'E24' '  -5.50%  '
Set-TextValue 'D25' '2.343'
Set-TextValue 'E25' '  -2.90%  '
Set-TextValue 'D26' '143.83'
Set-TextValue 'E26' '  -2.15%  '
Set-TextValue 'D27' '1.646'
Set-TextValue 'E27' '  -1.23%  '
Set-TextValue 'D28' '16.98'
Set-TextValue 'E28' '  -3.98%  '
Set-TextValue 'D29' '111.76'
Set-TextValue 'E29' '  -2.06%  '
Set-TextValue 'D30' '4.204'
Set-TextValue 'E30' '  -4.79%  '
Set-TextValue 'D31' '4.157'
Set-TextValue 'E31' '  -5.57%  '
Set-TextValue 'D32' '0.08796'
Set-TextValue 'E32' '  -1.03%  '
Set-TextValue 'D33' '0.04826'
Set-TextValue 'E33' '  -1.99%  '
Set-TextValue 'D34' '1.136'
Set-TextValue 'E34' '  -3.79%  '
Set-TextValue 'D35' '2.857'
Set-TextValue 'E35' '  -1.34%  '
Set-TextValue 'D37' '3.113'
Set-TextValue 'E37' '  -5.67%  '
Set-TextValue 'D38' '0.01795'
Set-TextValue 'E38' '  -5.40%  '
Set-TextValue 'D39' '2.221'
Set-TextValue 'E39' '  -7.80%  '
Set-TextValue 'D40' '0.4915'
Set-TextValue 'E40' '  -7.76%  '
Set-TextValue 'D41' '112.35'
Set-TextValue 'E41' '  -4.15%  '
Set-TextValue 'D42' '0.8989'
Set-TextValue 'E42' '  -8.51%  '
Set-TextValue 'D43' '6.186'
Set-TextValue 'E43' '  -2.22%  '
Set-TextValue 'D44' '1.001'
Set-TextValue 'E44' '  +0.08%  '
Set-TextValue 'D45' '7.786'
Set-TextValue 'E45' '  -5.24%  '
Set-TextValue 'D46' '0.4199'
Set-TextValue 'E46' '  -9.55%  '
Set-TextValue 'D47' '0.1267'
Set-TextValue 'E47' '  -7.71%  '
Set-TextValue 'D48' '9.139'
Set-TextValue 'E48' '  -3.76%  '
Set-TextValue 'D49' '0.05883'
Set-TextValue 'E49' '  -1.14%  '
Set-TextValue 'D50' '35.34'
Set-TextValue 'E50' '  -4.24%  '
Set-TextValue 'B51' 'NEARProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D51' '1.416'
Set-TextValue 'E51' '  -7.49%  '
